$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The R script appended two new rows of OHLC data (rows 152 and 153) below
# the existing history. Replicate the formatting of the preceding rows
# (date column keeps its yyyy-mm-dd hh:mm:ss style, G/H stay text) by
# copying format/values from existing rows instead of creating new styles.

# Copy the date-column number format (style s="1") down into the new rows.
$ws.Range("A151").Copy()
$ws.Range("A152:A153").PasteSpecial(-4122)  # xlPasteFormats

# Row 152
$ws.Range("A152").Value = 45454.2916666667
$ws.Range("B152").Value = 0
$ws.Range("C152").Value = 2.88000011444092
$ws.Range("D152").Value = 2.88000011444092
$ws.Range("E152").Value = 2.88000011444092
$ws.Range("F152").Value = 2.88000011444092

# G152/H152 reuse the same text values already present in row 151, so copy
# them as values to preserve their text (shared-string) type.
$ws.Range("G151").Copy()
$ws.Range("G152").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("H151").Copy()
$ws.Range("H152").PasteSpecial(-4163)  # xlPasteValues

# Row 153
$ws.Range("A153").Value = 45455.6119791667
$ws.Range("B153").Value = 3000
$ws.Range("C153").Value = 2.98000001907349
$ws.Range("D153").Value = 2.83999991416931
$ws.Range("E153").Value = 2.88000011444092
$ws.Range("F153").Value = 2.96000003814697

# G153 needs the text value "2.96000003814697", which already exists as a
# shared string elsewhere in the sheet (e.g. G57); H153 is "XHS.MI" (H151).
$ws.Range("G57").Copy()
$ws.Range("G153").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("H151").Copy()
$ws.Range("H153").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
